# xlsx,csv: trim whitespace, yield to optional block
# Closes #108 and #79
#
# On the "Resources" sheet: the Wood/Metal/Stone row labels pick up the
# (previously trimmed) trailing whitespace straight from the source data,
# and a new "Qty" column is added alongside the existing "Cost" column.

$wb = $excel.ActiveWorkbook

$wsExpansion = $wb.Worksheets.Item("Expansion")
$wsResources = $wb.Worksheets.Item("Resources")

# --- Resources sheet: add the new Qty column next to Cost, then pad the
#     material names with their untrimmed whitespace. ---
$wsResources.Range("C1").Value = "Qty"
$wsResources.Range("C1").Font.Bold = $true

$wsResources.Range("C2").Value = 3
$wsResources.Range("C3").Value = 2
$wsResources.Range("C4").Value = 1

$wsResources.Range("A2").Value = "Wood                                                   "
$wsResources.Range("A3").Value = "Metal                                                 "
$wsResources.Range("A4").Value = "Stone                           "

# --- Cursor / selection bookkeeping -----------------------------------
# Expansion's remembered selection moves up one row (B4 -> B3); restore
# the originally active sheet (Resources) afterwards so tab selection
# isn't disturbed.
[void]$wsExpansion.Range("B3").Select()
[void]$wsResources.Select()
[void]$wsResources.Range("N9").Select()
